$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the date column (A2:A6) from 2025-12-05 to 2025-12-06 ---
# Force the cells to stay as text (not auto-converted to a date serial)
# by pre-setting a text number format, then reset the style back to
# Normal so no explicit style index is left on the cells.
$dateRange = $ws.Range("A2:A6")
$dateRange.NumberFormat = "@"
$dateRange.Value = "2025-12-06"
$dateRange.Style = "Normal"

# --- Update K column (score) values ---
$ws.Range("K2").Value = 54.9
$ws.Range("K3").Value = 50.9
$ws.Range("K4").Value = 46.7
$ws.Range("K5").Value = 45.5
$ws.Range("K6").Value = 36.9

# --- Update N column (MACRO_SCORE) values ---
$ws.Range("N2").Value = 51.54219175917372
$ws.Range("N3").Value = 51.54219175917372
$ws.Range("N4").Value = 51.54219175917372
$ws.Range("N5").Value = 51.54219175917372
$ws.Range("N6").Value = 51.54219175917372
